$d = $word.ActiveDocument

# Remove the incorrect "beta ratio" statement: the resume originally read
# "...company's P/E and Beta Ratio, Income Statement..." which should read
# "...company's P/E ratio, Beta, Income Statement..."
$find = [char]0x2019 + "s P/E and Beta Ratio, Income Statement"
$replace = [char]0x2019 + "s P/E ratio, Beta, Income Statement"

$d.Content.Find.Execute($find, $true, $false, $false, $false, $false, `
                         $true, 1, $false, $replace, 2)
